# Fixing network data cleaning scripts
# - Rename header columns to clean machine-readable names
# - Normalize "de/del/los/la" -> "De/Del/Los/La" in several place names
# - Fix a floating point rounding value in D101
# - Remove the trailing footnote/metadata rows (150-154) and shrink the
#   worksheet dimension accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames ---------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case "de/del/los/la" particles in place names ---------------
$ws.Range("B5").Value = "Ocozocoautla De Espinosa"
$ws.Range("B9").Value = "Hidalgo Del Parral"
$ws.Range("A16").Value = "Ciudad De México"
$ws.Range("A28").Value = "Estado De México"
$ws.Range("B29").Value = "Atizapán De Zaragoza"
$ws.Range("B33").Value = "Tlalnepantla De Baz"
$ws.Range("B38").Value = "San Miguel De Allende"
$ws.Range("B39").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B44").Value = "Valle De Santiago"
$ws.Range("B47").Value = "Acapulco De Juárez"
$ws.Range("B48").Value = "Chilpancingo De Los Bravo"
$ws.Range("B53").Value = "Santiago De Anaya"
$ws.Range("B54").Value = "Tepehuacán De Guerrero"
$ws.Range("B58").Value = "Autlán De Navarro"
$ws.Range("B66").Value = "La Manzanilla De La Paz"
$ws.Range("B67").Value = "Ojuelos De Jalisco"
$ws.Range("B69").Value = "San Juan De Los Lagos"
$ws.Range("B70").Value = "Talpa De Allende"
$ws.Range("B72").Value = "Tepatitlán De Morelos"
$ws.Range("B77").Value = "Zapotlán Del Rey"
$ws.Range("B94").Value = "Ixtlán Del Río"
$ws.Range("B102").Value = "Chalcatongo De Hidalgo"
$ws.Range("B103").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B104").Value = "Oaxaca De Juárez"
$ws.Range("B112").Value = "Teotitlán De Flores Magón"
$ws.Range("B113").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B120").Value = "Cuayuca De Andrade"
$ws.Range("B122").Value = "Los Reyes De Juárez"
$ws.Range("B128").Value = "Cadereyta De Montes"
$ws.Range("B131").Value = "Mexquitic De Carmona"
$ws.Range("B140").Value = "Hueyapan De Ocampo"
$ws.Range("B141").Value = "Ignacio De La Llave"

# --- Floating point fix -------------------------------------------------
$ws.Range("D101").Value = 0.09595959595959597

# --- Remove trailing footnote / metadata rows (150-154) ----------------
$ws.Range("A150:A154").EntireRow.Delete()
